# "Generate Report for Handoff"
#
# This script moves the localization-status report from the
# "In Translation" state to "Ready for handoff": it updates the status
# text and the related timestamps on all three worksheets, and widens
# the Status columns that now have to fit the longer "Ready for handoff"
# label.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Status: "In Translation" -> "Ready for handoff"
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value     = "Ready for handoff"   # Status
$wsDeDe.Range("C2").Value     = "Ready for handoff"   # Status

# ---------------------------------------------------------------------
# Timestamps refreshed as part of the handoff-file regeneration
# ---------------------------------------------------------------------
# Overview!G2 ("Latest HO Xliff Generate Date") and de-de!H2
# ("Latest Handoff Datetime") shared the same original value.
$wsOverview.Range("G2").Value = "2016-08-29 23:02:12"
$wsDeDe.Range("H2").Value     = "2016-08-29 23:02:12"

# zh-cn!H2 ("Latest Handoff Datetime") had its own timestamp.
$wsZhCn.Range("H2").Value = "2016-08-29 23:02:00"

# ---------------------------------------------------------------------
# Widen the Status columns so the new "Ready for handoff" text fits
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333   # zh-cn column
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333   # de-de column
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.3333333333333   # Status column
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.3333333333333   # Status column
